$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.359656
$ws.Range("H2").Value = 7.078968
$ws.Range("I2").Value = 0.1135804410355361
$ws.Range("J2").Value = 0.1135804410355361
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8109183333333333
$ws.Range("N2").Value = 2.432755
$ws.Range("O2").Value = 0.2370884169621149
$ws.Range("P2").Value = 0.2370884169621149
$ws.Range("Q2").Value = 1.91348831076
$ws.Range("R2").Value = 17.22139479684
$ws.Range("S2").Value = 0.02692860696297409
$ws.Range("T2").Value = 0.0269286069629741

$ws.Range("G3").Value = 2.359656
$ws.Range("H3").Value = 7.078968
$ws.Range("I3").Value = 0.1135804410355361
$ws.Range("J3").Value = 0.1135804410355361
$ws.Range("O3").Value = 0.3378801459239538
$ws.Range("P3").Value = 0.3378801459239539
$ws.Range("Q3").Value = 2.7269561202
$ws.Range("R3").Value = 24.5426050818
$ws.Range("S3").Value = 0.03837657599119398
$ws.Range("T3").Value = 0.03837657599119398

$ws.Range("G4").Value = 2.359656
$ws.Range("H4").Value = 7.078968
$ws.Range("I4").Value = 0.1135804410355361
$ws.Range("J4").Value = 0.1135804410355361
$ws.Range("M4").Value = 1.434534666666667
$ws.Range("N4").Value = 4.303604
$ws.Range("O4").Value = 0.4194152964814894
$ws.Range("P4").Value = 0.4194152964814894
$ws.Range("Q4").Value = 3.385008333408
$ws.Range("R4").Value = 30.465075000672
$ws.Range("S4").Value = 0.0476373743514177
$ws.Range("T4").Value = 0.0476373743514177

$ws.Range("G5").Value = 2.359656
$ws.Range("H5").Value = 7.078968
$ws.Range("I5").Value = 0.1135804410355361
$ws.Range("J5").Value = 0.1135804410355361
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.019209
$ws.Range("N5").Value = 0.057627
$ws.Range("O5").Value = 0.005616140632441737
$ws.Range("P5").Value = 0.005616140632441737
$ws.Range("Q5").Value = 0.04532663210399999
$ws.Range("R5").Value = 0.407939688936
$ws.Range("S5").Value = 0.0006378837299503272
$ws.Range("T5").Value = 0.0006378837299503273

$ws.Range("I6").Value = 0.07630393871923234
$ws.Range("J6").Value = 0.07630393871923234
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8109183333333333
$ws.Range("N6").Value = 2.432755
$ws.Range("O6").Value = 0.2370884169621149
$ws.Range("P6").Value = 0.2370884169621149
$ws.Range("Q6").Value = 1.285491528937778
$ws.Range("R6").Value = 11.56942376044
$ws.Range("S6").Value = 0.01809078003891702
$ws.Range("T6").Value = 0.01809078003891702

$ws.Range("I7").Value = 0.07630393871923234
$ws.Range("J7").Value = 0.07630393871923234
$ws.Range("O7").Value = 0.3378801459239538
$ws.Range("P7").Value = 0.3378801459239539
$ws.Range("S7").Value = 0.02578158594902666
$ws.Range("T7").Value = 0.02578158594902666

$ws.Range("I8").Value = 0.07630393871923234
$ws.Range("J8").Value = 0.07630393871923234
$ws.Range("M8").Value = 1.434534666666667
$ws.Range("N8").Value = 4.303604
$ws.Range("O8").Value = 0.4194152964814894
$ws.Range("P8").Value = 0.4194152964814894
$ws.Range("Q8").Value = 2.274066433283556
$ws.Range("R8").Value = 20.466597899552
$ws.Range("S8").Value = 0.03200303908063223
$ws.Range("T8").Value = 0.03200303908063223

$ws.Range("I9").Value = 0.07630393871923234
$ws.Range("J9").Value = 0.07630393871923234
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.019209
$ws.Range("N9").Value = 0.057627
$ws.Range("O9").Value = 0.005616140632441737
$ws.Range("P9").Value = 0.005616140632441737
$ws.Range("Q9").Value = 0.030450670264
$ws.Range("R9").Value = 0.274056032376
$ws.Range("S9").Value = 0.000428533650656425
$ws.Range("T9").Value = 0.0004285336506564251

$ws.Range("G10").Value = 3.018243333333333
$ws.Range("H10").Value = 9.05473
$ws.Range("I10").Value = 0.1452810956141771
$ws.Range("J10").Value = 0.1452810956141771
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8109183333333333
$ws.Range("N10").Value = 2.432755
$ws.Range("O10").Value = 0.2370884169621149
$ws.Range("P10").Value = 0.2370884169621149
$ws.Range("Q10").Value = 2.447548853461111
$ws.Range("R10").Value = 22.02793968115
$ws.Range("S10").Value = 0.0344444649736869
$ws.Range("T10").Value = 0.03444446497368691

$ws.Range("G11").Value = 3.018243333333333
$ws.Range("H11").Value = 9.05473
$ws.Range("I11").Value = 0.1452810956141771
$ws.Range("J11").Value = 0.1452810956141771
$ws.Range("O11").Value = 0.3378801459239538
$ws.Range("P11").Value = 0.3378801459239539
$ws.Range("Q11").Value = 3.488058060194444
$ws.Range("R11").Value = 31.39252254174999
$ws.Range("S11").Value = 0.04908759778611005
$ws.Range("T11").Value = 0.04908759778611006

$ws.Range("G12").Value = 3.018243333333333
$ws.Range("H12").Value = 9.05473
$ws.Range("I12").Value = 0.1452810956141771
$ws.Range("J12").Value = 0.1452810956141771
$ws.Range("M12").Value = 1.434534666666667
$ws.Range("N12").Value = 4.303604
$ws.Range("O12").Value = 0.4194152964814894
$ws.Range("P12").Value = 0.4194152964814894
$ws.Range("Q12").Value = 4.329774694102221
$ws.Range("R12").Value = 38.96797224692
$ws.Range("S12").Value = 0.06093311379017569
$ws.Range("T12").Value = 0.06093311379017569

$ws.Range("G13").Value = 3.018243333333333
$ws.Range("H13").Value = 9.05473
$ws.Range("I13").Value = 0.1452810956141771
$ws.Range("J13").Value = 0.1452810956141771
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.019209
$ws.Range("N13").Value = 0.057627
$ws.Range("O13").Value = 0.005616140632441737
$ws.Range("P13").Value = 0.005616140632441737
$ws.Range("Q13").Value = 0.05797743619
$ws.Range("R13").Value = 0.5217969257099999
$ws.Range("S13").Value = 0.0008159190642044329
$ws.Range("T13").Value = 0.0008159190642044331

$ws.Range("G14").Value = 1.732629
$ws.Range("H14").Value = 5.197887
$ws.Range("I14").Value = 0.083398921694925
$ws.Range("J14").Value = 0.083398921694925
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.8109183333333333
$ws.Range("N14").Value = 2.432755
$ws.Range("O14").Value = 0.2370884169621149
$ws.Range("P14").Value = 0.2370884169621149
$ws.Range("Q14").Value = 1.405020620965
$ws.Range("R14").Value = 12.645185588685
$ws.Range("S14").Value = 0.01977291832099715
$ws.Range("T14").Value = 0.01977291832099715

$ws.Range("G15").Value = 1.732629
$ws.Range("H15").Value = 5.197887
$ws.Range("I15").Value = 0.083398921694925
$ws.Range("J15").Value = 0.083398921694925
$ws.Range("O15").Value = 0.3378801459239538
$ws.Range("P15").Value = 0.3378801459239539
$ws.Range("Q15").Value = 2.002327142425
$ws.Range("R15").Value = 18.020944281825
$ws.Range("S15").Value = 0.02817883983218165
$ws.Range("T15").Value = 0.02817883983218166

$ws.Range("G16").Value = 1.732629
$ws.Range("H16").Value = 5.197887
$ws.Range("I16").Value = 0.083398921694925
$ws.Range("J16").Value = 0.083398921694925
$ws.Range("M16").Value = 1.434534666666667
$ws.Range("N16").Value = 4.303604
$ws.Range("O16").Value = 0.4194152964814894
$ws.Range("P16").Value = 0.4194152964814894
$ws.Range("Q16").Value = 2.485516364972
$ws.Range("R16").Value = 22.369647284748
$ws.Range("S16").Value = 0.03497878346891348
$ws.Range("T16").Value = 0.03497878346891348

$ws.Range("G17").Value = 1.732629
$ws.Range("H17").Value = 5.197887
$ws.Range("I17").Value = 0.083398921694925
$ws.Range("J17").Value = 0.083398921694925
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.019209
$ws.Range("N17").Value = 0.057627
$ws.Range("O17").Value = 0.005616140632441737
$ws.Range("P17").Value = 0.005616140632441737
$ws.Range("Q17").Value = 0.033282070461
$ws.Range("R17").Value = 0.2995386341489999
$ws.Range("S17").Value = 0.0004683800728326949
$ws.Range("T17").Value = 0.000468380072832695

$ws.Range("G18").Value = 9.589644
$ws.Range("H18").Value = 28.768932
$ws.Range("I18").Value = 0.4615910094072114
$ws.Range("J18").Value = 0.4615910094072114
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.8109183333333333
$ws.Range("N18").Value = 2.432755
$ws.Range("O18").Value = 0.2370884169621149
$ws.Range("P18").Value = 0.2370884169621149
$ws.Range("Q18").Value = 7.77641812974
$ws.Range("R18").Value = 69.98776316765999
$ws.Range("S18").Value = 0.1094378817043004
$ws.Range("T18").Value = 0.1094378817043005

$ws.Range("G19").Value = 9.589644
$ws.Range("H19").Value = 28.768932
$ws.Range("I19").Value = 0.4615910094072114
$ws.Range("J19").Value = 0.4615910094072114
$ws.Range("O19").Value = 0.3378801459239538
$ws.Range("P19").Value = 0.3378801459239539
$ws.Range("Q19").Value = 11.0823520023
$ws.Range("R19").Value = 99.74116802069999
$ws.Range("S19").Value = 0.1559624376156937
$ws.Range("T19").Value = 0.1559624376156938

$ws.Range("G20").Value = 9.589644
$ws.Range("H20").Value = 28.768932
$ws.Range("I20").Value = 0.4615910094072114
$ws.Range("J20").Value = 0.4615910094072114
$ws.Range("M20").Value = 1.434534666666667
$ws.Range("N20").Value = 4.303604
$ws.Range("O20").Value = 0.4194152964814894
$ws.Range("P20").Value = 0.4194152964814894
$ws.Range("Q20").Value = 13.756676758992
$ws.Range("R20").Value = 123.810090830928
$ws.Range("S20").Value = 0.1935983300637155
$ws.Range("T20").Value = 0.1935983300637155

$ws.Range("G21").Value = 9.589644
$ws.Range("H21").Value = 28.768932
$ws.Range("I21").Value = 0.4615910094072114
$ws.Range("J21").Value = 0.4615910094072114
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0.3333333333333333
$ws.Range("M21").Value = 0.019209
$ws.Range("N21").Value = 0.057627
$ws.Range("O21").Value = 0.005616140632441737
$ws.Range("P21").Value = 0.005616140632441737
$ws.Range("Q21").Value = 0.184207471596
$ws.Range("R21").Value = 1.657867244364
$ws.Range("S21").Value = 0.002592360023501636
$ws.Range("T21").Value = 0.002592360023501636

$ws.Range("G22").Value = 2.489795
$ws.Range("H22").Value = 7.469385
$ws.Range("I22").Value = 0.1198445935289181
$ws.Range("J22").Value = 0.1198445935289181
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 0.8109183333333333
$ws.Range("N22").Value = 2.432755
$ws.Range("O22").Value = 0.2370884169621149
$ws.Range("P22").Value = 0.2370884169621149
$ws.Range("Q22").Value = 2.019020411741666
$ws.Range("R22").Value = 18.171183705675
$ws.Range("S22").Value = 0.0284137649612393
$ws.Range("T22").Value = 0.02841376496123931

$ws.Range("G23").Value = 2.489795
$ws.Range("H23").Value = 7.469385
$ws.Range("I23").Value = 0.1198445935289181
$ws.Range("J23").Value = 0.1198445935289181
$ws.Range("O23").Value = 0.3378801459239538
$ws.Range("P23").Value = 0.3378801459239539
$ws.Range("Q23").Value = 2.877352340041667
$ws.Range("R23").Value = 25.896171060375
$ws.Range("S23").Value = 0.04049310874974776
$ws.Range("T23").Value = 0.04049310874974777

$ws.Range("G24").Value = 2.489795
$ws.Range("H24").Value = 7.469385
$ws.Range("I24").Value = 0.1198445935289181
$ws.Range("J24").Value = 0.1198445935289181
$ws.Range("M24").Value = 1.434534666666667
$ws.Range("N24").Value = 4.303604
$ws.Range("O24").Value = 0.4194152964814894
$ws.Range("P24").Value = 0.4194152964814894
$ws.Range("Q24").Value = 3.571697240393334
$ws.Range("R24").Value = 32.14527516354
$ws.Range("S24").Value = 0.05026465572663476
$ws.Range("T24").Value = 0.05026465572663476

$ws.Range("G25").Value = 2.489795
$ws.Range("H25").Value = 7.469385
$ws.Range("I25").Value = 0.1198445935289181
$ws.Range("J25").Value = 0.1198445935289181
$ws.Range("K25").Value = 1
$ws.Range("L25").Value = 0.3333333333333333
$ws.Range("M25").Value = 0.019209
$ws.Range("N25").Value = 0.057627
$ws.Range("O25").Value = 0.005616140632441737
$ws.Range("P25").Value = 0.005616140632441737
$ws.Range("Q25").Value = 0.047826472155
$ws.Range("R25").Value = 0.430438249395
$ws.Range("S25").Value = 0.0006730640912962207
$ws.Range("T25").Value = 0.0006730640912962208
